# Refresh the crypto price/volume table to the latest scrape.
# Cells in column D that look like plain numbers (e.g. "1.006") are written
# with a leading apostrophe so Excel keeps them as text (matching the sheet's
# existing "27.944.91" / "0.4818"-style price strings) instead of parsing them
# into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = '27.944.91'
    "E2" = '  +1.77%  '
    "D3" = '1.902.91'
    "E3" = '  +2.35%  '
    "D4" = '''1.006'
    "E4" = '  -0.54%  '
    "D5" = '''316.89'
    "E5" = '  +0.42%  '
    "D6" = '''1.004'
    "E6" = '  -0.65%  '
    "D7" = '''0.4818'
    "E7" = '  +0.96%  '
    "D8" = '''0.3791'
    "E8" = '  -0.24%  '
    "D9" = '''0.07365'
    "E9" = '  +0.71%  '
    "D10" = '''0.9313'
    "E10" = '  +0.14%  '
    "D11" = '''20.74'
    "E11" = '  +0.10%  '
    "B12" = 'TRON'
    "C12" = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    "D12" = '''0.07737'
    "E12" = '  -0.71%  '
    "B13" = 'WrappedEther'
    "C13" = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    "D13" = '1.948.09'
    "E13" = '  +4.43%  '
    "D14" = '''5.475'
    "E14" = '  +0.62%  '
    "D15" = '''6.629'
    "E15" = '  +1.25%  '
    "D16" = '''91.68'
    "E16" = '  +1.63%  '
    "D17" = '''1.004'
    "E17" = '  -0.67%  '
    "D18" = '''0.000008860'
    "E18" = '  +0.48%  '
    "E19" = '  -0.62%  '
    "D20" = '28.010.20'
    "E20" = '  +1.81%  '
    "D21" = '''14.66'
    "E21" = '  +0.23%  '
    "D22" = '''5.148'
    "E22" = '  +1.02%  '
    "D23" = '2.188.15'
    "E23" = '  +3.95%  '
    "D24" = '''10.89'
    "E24" = '  +1.97%  '
    "D25" = '''156.12'
    "E25" = '  +0.81%  '
    "D26" = '''1.909'
    "E26" = '  -1.83%  '
    "D27" = '''18.44'
    "E27" = '  +0.01%  '
    "D28" = '''2.119'
    "E28" = '  +5.81%  '
    "D29" = '''117.06'
    "E29" = '  +1.51%  '
    "D30" = '''4.965'
    "E30" = '  +0.49%  '
    "D31" = '''0.08928'
    "E31" = '  +0.38%  '
    "D32" = '''3.260'
    "E32" = '  -2.17%  '
    "D33" = '''1.250'
    "E33" = '  +3.76%  '
    "D34" = '''0.7679'
    "E34" = '  +2.00%  '
    "D35" = '''4.664'
    "E35" = '  +1.85%  '
    "D36" = '''2.584'
    "E36" = '  -4.25%  '
    "D37" = '''0.02055'
    "E37" = '  +0.52%  '
    "D38" = '''1.104'
    "E38" = '  -1.84%  '
    "D39" = '''0.5487'
    "E39" = '  -1.34%  '
    "D40" = '''3.001'
    "E40" = '  +0.47%  '
    "D41" = '''0.05268'
    "E41" = '  -0.14%  '
    "D42" = '''6.926'
    "E42" = '  -1.38%  '
    "D43" = '''0.1524'
    "E43" = '  +0.64%  '
    "D44" = '''8.481'
    "E44" = '  -1.03%  '
    "D45" = '''110.09'
    "E45" = '  +6.73%  '
    "D46" = '''10.68'
    "E46" = '  +0.51%  '
    "D47" = '''0.4810'
    "E47" = '  -1.21%  '
    "D48" = '''1.004'
    "E48" = '  -0.72%  '
    "E49" = '  -0.95%  '
    "D50" = '''67.97'
    "E50" = '  +0.84%  '
    "E51" = '  -0.31%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
